$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 12,20
$arr[0,0] = "ECs"
$arr[0,1] = "Fgf1"
$arr[0,2] = "Fgfr2"
$arr[0,3] = "ECs"
$arr[0,4] = [double]"3"
$arr[0,5] = [double]"1"
$arr[0,6] = [double]"1.384145666666667"
$arr[0,7] = [double]"4.152437"
$arr[0,8] = [double]"0.1014617184198512"
$arr[0,9] = [double]"0.1334061399754118"
$arr[0,10] = [double]"2"
$arr[0,11] = [double]"0.6666666666666666"
$arr[0,12] = [double]"0.493831"
$arr[0,13] = [double]"1.481493"
$arr[0,14] = [double]"0.1121895146192186"
$arr[0,15] = [double]"0.1134277760249069"
$arr[0,16] = [double]"0.6835340387156666"
$arr[0,17] = [double]"6.151806348440999"
$arr[0,18] = [double]"0.01138294094195493"
$arr[0,19] = [double]"0.01513196176547839"
$arr[1,0] = "ECs"
$arr[1,1] = "Fgf1"
$arr[1,2] = "Fgfr2"
$arr[1,3] = "FAPs"
$arr[1,4] = [double]"3"
$arr[1,5] = [double]"1"
$arr[1,6] = [double]"1.384145666666667"
$arr[1,7] = [double]"4.152437"
$arr[1,8] = [double]"0.1014617184198512"
$arr[1,9] = [double]"0.1334061399754118"
$arr[1,10] = [double]"3"
$arr[1,11] = [double]"1"
$arr[1,12] = [double]"3.763360333333333"
$arr[1,13] = [double]"11.290081"
$arr[1,14] = [double]"0.85496773012202"
$arr[1,15] = [double]"0.8644042050627692"
$arr[1,16] = [double]"5.209038897488556"
$arr[1,17] = [double]"46.881350077397"
$arr[1,18] = [double]"0.0867464950916997"
$arr[1,19] = [double]"0.1153168283759383"
$arr[2,0] = "ECs"
$arr[2,1] = "Fgf1"
$arr[2,2] = "Fgfr2"
$arr[2,3] = "M1"
$arr[2,4] = [double]"3"
$arr[2,5] = [double]"1"
$arr[2,6] = [double]"1.384145666666667"
$arr[2,7] = [double]"4.152437"
$arr[2,8] = [double]"0.1014617184198512"
$arr[2,9] = [double]"0.1334061399754118"
$arr[2,10] = [double]"1"
$arr[2,11] = [double]"0.3333333333333333"
$arr[2,12] = [double]"0.0004073333333333333"
$arr[2,13] = [double]"0.001222"
$arr[2,14] = [double]"9.253880164447964E-05"
$arr[2,15] = [double]"9.356017362379453E-05"
$arr[2,16] = [double]"0.0005638086682222222"
$arr[2,17] = [double]"0.005074278014"
$arr[2,18] = [double]"9.389145835362654E-06"
$arr[2,19] = [double]"1.248150161857976E-05"
$arr[3,0] = "ECs"
$arr[3,1] = "Fgf1"
$arr[3,2] = "Fgfr2"
$arr[3,3] = "sCs"
$arr[3,4] = [double]"3"
$arr[3,5] = [double]"1"
$arr[3,6] = [double]"1.384145666666667"
$arr[3,7] = [double]"4.152437"
$arr[3,8] = [double]"0.1014617184198512"
$arr[3,9] = [double]"0.1334061399754118"
$arr[3,10] = [double]"2"
$arr[3,11] = [double]"1"
$arr[3,12] = [double]"0.1441585"
$arr[3,13] = [double]"0.288317"
$arr[3,14] = [double]"0.03275021645711715"
$arr[3,15] = [double]"0.02207445873870014"
$arr[3,16] = [double]"0.1995363630881667"
$arr[3,17] = [double]"1.197218178529"
$arr[3,18] = [double]"0.003322893240361196"
$arr[3,19] = [double]"0.002944868332376482"
$arr[4,0] = "FAPs"
$arr[4,1] = "Fgf1"
$arr[4,2] = "Fgfr2"
$arr[4,3] = "ECs"
$arr[4,4] = [double]"3"
$arr[4,5] = [double]"1"
$arr[4,6] = [double]"2.458038666666667"
$arr[4,7] = [double]"7.374116000000001"
$arr[4,8] = [double]"0.1801810554109116"
$arr[4,9] = [double]"0.2369096391566985"
$arr[4,10] = [double]"2"
$arr[4,11] = [double]"0.6666666666666666"
$arr[4,12] = [double]"0.493831"
$arr[4,13] = [double]"1.481493"
$arr[4,14] = [double]"0.1121895146192186"
$arr[4,15] = [double]"0.1134277760249069"
$arr[4,16] = [double]"1.213855692798667"
$arr[4,17] = [double]"10.924701235188"
$arr[4,18] = [double]"0.02021442515012869"
$arr[4,19] = [double]"0.02687213348840752"
$arr[5,0] = "FAPs"
$arr[5,1] = "Fgf1"
$arr[5,2] = "Fgfr2"
$arr[5,3] = "FAPs"
$arr[5,4] = [double]"3"
$arr[5,5] = [double]"1"
$arr[5,6] = [double]"2.458038666666667"
$arr[5,7] = [double]"7.374116000000001"
$arr[5,8] = [double]"0.1801810554109116"
$arr[5,9] = [double]"0.2369096391566985"
$arr[5,10] = [double]"3"
$arr[5,11] = [double]"1"
$arr[5,12] = [double]"3.763360333333333"
$arr[5,13] = [double]"11.290081"
$arr[5,14] = [double]"0.85496773012202"
$arr[5,15] = [double]"0.8644042050627692"
$arr[5,16] = [double]"9.250485215932891"
$arr[5,17] = [double]"83.25436694339601"
$arr[5,18] = [double]"0.154048987955657"
$arr[5,19] = [double]"0.2047856883069535"
$arr[6,0] = "FAPs"
$arr[6,1] = "Fgf1"
$arr[6,2] = "Fgfr2"
$arr[6,3] = "M1"
$arr[6,4] = [double]"3"
$arr[6,5] = [double]"1"
$arr[6,6] = [double]"2.458038666666667"
$arr[6,7] = [double]"7.374116000000001"
$arr[6,8] = [double]"0.1801810554109116"
$arr[6,9] = [double]"0.2369096391566985"
$arr[6,10] = [double]"1"
$arr[6,11] = [double]"0.3333333333333333"
$arr[6,12] = [double]"0.0004073333333333333"
$arr[6,13] = [double]"0.001222"
$arr[6,14] = [double]"9.253880164447964E-05"
$arr[6,15] = [double]"9.356017362379453E-05"
$arr[6,16] = [double]"0.001001241083555556"
$arr[6,17] = [double]"0.009011169752000002"
$arr[6,18] = [double]"1.667373894676334E-05"
$arr[6,19] = [double]"2.216530697265123E-05"
$arr[7,0] = "FAPs"
$arr[7,1] = "Fgf1"
$arr[7,2] = "Fgfr2"
$arr[7,3] = "sCs"
$arr[7,4] = [double]"3"
$arr[7,5] = [double]"1"
$arr[7,6] = [double]"2.458038666666667"
$arr[7,7] = [double]"7.374116000000001"
$arr[7,8] = [double]"0.1801810554109116"
$arr[7,9] = [double]"0.2369096391566985"
$arr[7,10] = [double]"2"
$arr[7,11] = [double]"1"
$arr[7,12] = [double]"0.1441585"
$arr[7,13] = [double]"0.288317"
$arr[7,14] = [double]"0.03275021645711715"
$arr[7,15] = [double]"0.02207445873870014"
$arr[7,16] = [double]"0.3543471671286667"
$arr[7,17] = [double]"2.126083002772"
$arr[7,18] = [double]"0.005900968566179172"
$arr[7,19] = [double]"0.00522965205436488"
$arr[8,0] = "sCs"
$arr[8,1] = "Fgf1"
$arr[8,2] = "Fgfr2"
$arr[8,3] = "ECs"
$arr[8,4] = [double]"2"
$arr[8,5] = [double]"1"
$arr[8,6] = [double]"9.799863999999999"
$arr[8,7] = [double]"19.599728"
$arr[8,8] = [double]"0.7183572261692373"
$arr[8,9] = [double]"0.6296842208678898"
$arr[8,10] = [double]"2"
$arr[8,11] = [double]"0.6666666666666666"
$arr[8,12] = [double]"0.493831"
$arr[8,13] = [double]"1.481493"
$arr[8,14] = [double]"0.1121895146192186"
$arr[8,15] = [double]"0.1134277760249069"
$arr[8,16] = [double]"4.839476638983999"
$arr[8,17] = [double]"29.036859833904"
$arr[8,18] = [double]"0.08059214852713494"
$arr[8,19] = [double]"0.07142368077102101"
$arr[9,0] = "sCs"
$arr[9,1] = "Fgf1"
$arr[9,2] = "Fgfr2"
$arr[9,3] = "FAPs"
$arr[9,4] = [double]"2"
$arr[9,5] = [double]"1"
$arr[9,6] = [double]"9.799863999999999"
$arr[9,7] = [double]"19.599728"
$arr[9,8] = [double]"0.7183572261692373"
$arr[9,9] = [double]"0.6296842208678898"
$arr[9,10] = [double]"3"
$arr[9,11] = [double]"1"
$arr[9,12] = [double]"3.763360333333333"
$arr[9,13] = [double]"11.290081"
$arr[9,14] = [double]"0.85496773012202"
$arr[9,15] = [double]"0.8644042050627692"
$arr[9,16] = [double]"36.88041944966133"
$arr[9,17] = [double]"221.282516697968"
$arr[9,18] = [double]"0.6141722470746633"
$arr[9,19] = [double]"0.5443016883798775"
$arr[10,0] = "sCs"
$arr[10,1] = "Fgf1"
$arr[10,2] = "Fgfr2"
$arr[10,3] = "M1"
$arr[10,4] = [double]"2"
$arr[10,5] = [double]"1"
$arr[10,6] = [double]"9.799863999999999"
$arr[10,7] = [double]"19.599728"
$arr[10,8] = [double]"0.7183572261692373"
$arr[10,9] = [double]"0.6296842208678898"
$arr[10,10] = [double]"1"
$arr[10,11] = [double]"0.3333333333333333"
$arr[10,12] = [double]"0.0004073333333333333"
$arr[10,13] = [double]"0.001222"
$arr[10,14] = [double]"9.253880164447964E-05"
$arr[10,15] = [double]"9.356017362379453E-05"
$arr[10,16] = [double]"0.003991811269333333"
$arr[10,17] = [double]"0.023950867616"
$arr[10,18] = [double]"6.647591686235365E-05"
$arr[10,19] = [double]"5.891336503256355E-05"
$arr[11,0] = "sCs"
$arr[11,1] = "Fgf1"
$arr[11,2] = "Fgfr2"
$arr[11,3] = "sCs"
$arr[11,4] = [double]"2"
$arr[11,5] = [double]"1"
$arr[11,6] = [double]"9.799863999999999"
$arr[11,7] = [double]"19.599728"
$arr[11,8] = [double]"0.7183572261692373"
$arr[11,9] = [double]"0.6296842208678898"
$arr[11,10] = [double]"2"
$arr[11,11] = [double]"1"
$arr[11,12] = [double]"0.1441585"
$arr[11,13] = [double]"0.288317"
$arr[11,14] = [double]"0.03275021645711715"
$arr[11,15] = [double]"0.02207445873870014"
$arr[11,16] = [double]"1.412733694444"
$arr[11,17] = [double]"5.650934777776"
$arr[11,18] = [double]"0.02352635465057678"
$arr[11,19] = [double]"0.01389993835195878"

$ws.Range("A2:T13").Value = $arr
